$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that were repulled/recalculated
$ws.Range("F8").Value = 11
$ws.Range("F10").Value = -2
$ws.Range("F12").Value = -9
$ws.Range("F20").Value = -1
$ws.Range("F23").Value = 9
$ws.Range("F26").Value = -9
$ws.Range("F32").Value = 1
